$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Fiscal Year From" row (row 6) to host the new
# "Charge Type" filter field. This shifts every row from 6 downward by one,
# which is exactly what the target layout needs (Fiscal Year From/To, Filter
# by, Period From/To, Date From/To, Run By, Run Date, the spacer row and the
# column-header row all move down one row without any other change).
$ws.Rows.Item(6).Insert()

# The freshly inserted row only carries formatting where the row above had
# content (columns A and B). Copy the full formatting of the row above
# (row 5, "Partner") across every used column so the new row matches the
# other label/input rows exactly (labels in A, blank input box in B, and
# plain cells beyond that), then set the new label text.
$ws.Range("A5:AX5").Copy()
$ws.Range("A6:AX6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A6").Value = "Charge Type"
